$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: Wins, Losses, Ties in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the same bold/centered/bordered header style used by the rest of row 1
# (e.g. AC1) by copying its formatting over rather than reconstructing it
# property-by-property.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Data rows 2-36: team win/loss/tie record, constant across all players on the roster
for ($r = 2; $r -le 36; $r++) {
    $ws.Cells.Item($r, 30).Value = 55
    $ws.Cells.Item($r, 31).Value = 60
    $ws.Cells.Item($r, 32).Value = 0
}
